$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "...the name of BrantEntry.bat. Then exit the script."
#        -> "...the name of BrantEntry.bat. Then exit and save [bookmark]the script."
#
# The single run containing "the name of BrantEntry.bat. Then exit the script."
# is split into three runs ("...Then exit ", "and save ", "the script."),
# with the (relocated) _GoBack bookmark sitting between the 2nd and 3rd.
# ------------------------------------------------------------------

$rng = $d.Content
$found = $rng.Find.Execute("the name of BrantEntry.bat. Then exit the script.")
if (-not $found) {
    throw "Could not find target sentence for change 1"
}

# Clear the matched range, then rebuild it piece by piece so each piece
# becomes its own run (InsertAfter on a freshly-collapsed range does not
# get coalesced into the preceding run).
$rng.Text = ""

$rng.InsertAfter("the name of BrantEntry.bat. Then exit ")
$piece2 = $d.Range($rng.End, $rng.End)
$piece2.InsertAfter("and save ")

$piece3Pos = $piece2.End
$piece3 = $d.Range($piece3Pos, $piece3Pos)
$piece3.InsertAfter("the script.")

# Re-seat the (single, special) _GoBack bookmark here; Word only ever keeps
# one _GoBack, so this automatically removes it from its old location.
$bmSpot = $d.Range($piece3Pos, $piece3Pos)
$d.Bookmarks.Add("_GoBack", $bmSpot)

# ------------------------------------------------------------------
# Change 2: "...NEST2020.csv; if [bookmark]you only write \\Brant-Data\\NEST2020.csv..."
#        -> "...NEST2020.csv; if you only write \\Brant-Data\\NEST2020.csv..."
#
# The old _GoBack bookmark that used to sit between ".csv; if " and
# "you only write " is gone now (moved above); merge those two runs back
# into a single run to match.
# ------------------------------------------------------------------

$rng2 = $d.Content
$found2 = $rng2.Find.Execute(".csv; if you only write ")
if (-not $found2) {
    throw "Could not find target phrase for change 2"
}
$innerStart = $rng2.Start
$innerEnd = $rng2.End

# Drop temporary barrier bookmarks just outside the two runs we want to
# merge, so the merge doesn't cascade into the neighboring runs
# ("\\Brant-Data\\Data\\NEST2020" on the left, "\\Brant-Data\\NEST2020.csv"
# on the right).
$leftBarrier = $d.Range($innerStart, $innerStart)
$d.Bookmarks.Add("TempLeft", $leftBarrier)
$rightBarrier = $d.Range($innerEnd, $innerEnd)
$d.Bookmarks.Add("TempRight", $rightBarrier)

$tl = $d.Bookmarks.Item("TempLeft")
$tr = $d.Bookmarks.Item("TempRight")
$mergeStart = $tl.Range.End
$mergeEnd = $tr.Range.Start

$mergeRng = $d.Range($mergeStart, $mergeEnd)
# Force an actual content change first (assigning the identical text is a
# no-op and would not trigger the run coalescing pass), then set the real
# text we want - this leaves a single merged run behind.
$mergeRng.Text = [char]1
$tl2 = $d.Bookmarks.Item("TempLeft")
$tr2 = $d.Bookmarks.Item("TempRight")
$mergeRng2 = $d.Range($tl2.Range.End, $tr2.Range.Start)
$mergeRng2.Text = ".csv; if you only write "

# Clean up the temporary barrier bookmarks.
$d.Bookmarks.Item("TempLeft").Delete()
$d.Bookmarks.Item("TempRight").Delete()
